$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text values (e.g. "211.00", "27.457.85").
# Force text format on those cells before assignment so Excel keeps them
# exactly as typed instead of converting to a Double and dropping
# formatting such as trailing zeros or thousand-grouping dots.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.457.85'
$ws.Range("E2").Value = '  -0.38%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.614.50'
$ws.Range("E3").Value = '  -1.57%  '

$ws.Range("E4").Value = '  +0.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.00'
$ws.Range("E5").Value = '  -0.86%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.525'
$ws.Range("E6").Value = '  -2.26%  '

$ws.Range("E7").Value = '  +0.26%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.90'
$ws.Range("E8").Value = '  -0.38%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.260'
$ws.Range("E9").Value = '  +1.11%  '

$ws.Range("E10").Value = '  -0.12%  '

$ws.Range("E11").Value = '  -0.45%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.845.66'
$ws.Range("E12").Value = '  -1.41%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.599.51'
$ws.Range("E13").Value = '  -2.62%  '

$ws.Range("E14").Value = '  -0.28%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.549'
$ws.Range("E15").Value = '  -2.59%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.23'
$ws.Range("E16").Value = '  +0.08%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.472.78'
$ws.Range("E17").Value = '  -0.09%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '227.08'
$ws.Range("E18").Value = '  -1.12%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0718'
$ws.Range("E19").Value = '  -0.78%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.52'
$ws.Range("E20").Value = '  -2.57%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.27'
$ws.Range("E22").Value = '  -0.88%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.88'
$ws.Range("E23").Value = '  -0.49%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.07'
$ws.Range("E24").Value = '  +6.49%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.96'
$ws.Range("E25").Value = '  -0.44%  '

$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.111'
$ws.Range("E26").Value = '  -1.55%  '

$ws.Range("B27").Value = 'BinanceUSD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.18%  '

$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.81'
$ws.Range("E28").Value = '  -2.36%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.52'
$ws.Range("E29").Value = '  -0.44%  '

$ws.Range("E30").Value = '  -0.80%  '

$ws.Range("E31").Value = '  -1.35%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.27'
$ws.Range("E32").Value = '  -0.50%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.442.65'
$ws.Range("E33").Value = '  +1.12%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.05'
$ws.Range("E34").Value = '  -3.63%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.52'
$ws.Range("E35").Value = '  -3.59%  '

$ws.Range("E36").Value = '  -0.08%  '

$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.943'
$ws.Range("E37").Value = '  +7.27%  '

$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.562'
$ws.Range("E38").Value = '  -1.57%  '

$ws.Range("E39").Value = '  +0.12%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.862'
$ws.Range("E40").Value = '  -2.01%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '69.06'
$ws.Range("E41").Value = '  +6.40%  '

$ws.Range("E42").Value = '  +0.25%  '

$ws.Range("E43").Value = '  -1.91%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.47'
$ws.Range("E44").Value = '  +0.23%  '

$ws.Range("E45").Value = '  -2.44%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.20'
$ws.Range("E46").Value = '  -2.18%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.756.64'
$ws.Range("E47").Value = '  -1.37%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.68'
$ws.Range("E48").Value = '  -0.03%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '86.06'
$ws.Range("E49").Value = '  -0.04%  '

$ws.Range("E50").Value = '  -0.87%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0986'
$ws.Range("E51").Value = '  -0.23%  '
